# Seq_Diag_Dependency_2: plot red arrows for inheritance and dependency (wip)
#
# The sequence-diagram grid keeps one "lifeline" column per nesting level in
# columns A..F (to the left of the "Class" header originally in G) plus the
# per-call-depth columns G..P ("Class" header / the numbered call-depth
# columns 0..8 and their arrows).
#
# This edit makes room for one more lifeline column (shifting the existing
# B:P data one column right, to C:Q - matching every "+1 column" move the
# diff shows), then:
#   - turns a handful of the existing "->" (single-line) call arrows into
#     the "=>" (double-line) dependency/inheritance arrow glyph
#   - adds a brand new "Farzi" dependency/inheritance arrow pair (B18 "=>"
#     lifeline glyph + G18 "▷" marker)
#   - extends the call-depth header row with two more levels ("9" / "10")
#     and the matching new arrows for the Farzi / check_farzi calls
#   - adds the new "Bike" inheritance lifeline glyph (B21 "◁")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new blank lifeline column before column B. ----------------
# This shifts the existing B:P columns (values, styles, and the column-width
# definitions) one to the right, to C:Q.
$ws.Columns("B:B").Insert()

# The insert picks up column A's thick-left-border style for every row of
# the new column B (Excel's usual "format from the cell to the left"
# behaviour). The diff only wants column B populated on rows 1, 18-21, so
# strip the unwanted carried-over formatting/values everywhere else first.
$ws.Range("B2:B17").Clear()

# --- 2. Promote a handful of "->" arrows to the "=>" dependency glyph. -----
$ws.Range("A2").Value = "⇒"
$ws.Range("C5").Value = "⇒"
$ws.Range("F5").Value = "⇒"
$ws.Range("E8").Value = "⇒"
$ws.Range("F8").Value = "⇒"
$ws.Range("D12").Value = "⇒"
$ws.Range("F12").Value = "⇒"

# --- 3. New "Farzi" dependency/inheritance arrow (row 18). -----------------
# B18 needs the thick-left-border lifeline styling used throughout columns
# A-F; copy formats from the already-styled A18 cell, then set its value.
$ws.Range("A18").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "⇒"
$ws.Range("G18").Value = "▷"

# --- 4. New trailing call-depth columns R ("9") and S ("10"). --------------
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R1").Value = "9"
$ws.Range("S1").PasteSpecial(-4122)
$ws.Range("S1").Value = "10"

# Style the new trailing cells in rows 19-24 like their row-neighbours
# (thick left border, no value) before filling in the actual arrows.
$ws.Range("Q19").Copy()
$ws.Range("R19").PasteSpecial(-4122)
$ws.Range("S19").PasteSpecial(-4122)
$ws.Range("R19").Value = "→"

$ws.Range("Q20").Copy()
$ws.Range("R20").PasteSpecial(-4122)
$ws.Range("S20").PasteSpecial(-4122)
$ws.Range("S20").Value = "→"

$ws.Range("Q21").Copy()
$ws.Range("R21").PasteSpecial(-4122)

$ws.Range("Q22").Copy()
$ws.Range("R22").PasteSpecial(-4122)

$ws.Range("Q23").Copy()
$ws.Range("R23").PasteSpecial(-4122)

$ws.Range("Q24").Copy()
$ws.Range("R24").PasteSpecial(-4122)

# --- 5. New lifeline arrow glyph for Bike's inheritance (row 21). ----------
$ws.Range("A21").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").Value = "◁"

# --- 6. Restore the original A1 selection (untouched by the diff). --------
$ws.Range("A1").Select()
